$d = $word.ActiveDocument

function ReplaceExact($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
    }
    return $ok
}

# 1. Merge "Creace & Querini" sentence into a single run (drop proofErr spell markers).
ReplaceExact "Tourism is now one of the largest industries and one of the fastest growing economic sectors around the globe (Creace & Querini, 2011). Furthermore, they noted that many countries had seen the potential of tourism for development because it drives new economic activity in a region." "Tourism is now one of the largest industries and one of the fastest growing economic sectors around the globe (Creace & Querini, 2011). Furthermore, they noted that many countries had seen the potential of tourism for development because it drives new economic activity in a region."

# 2. Merge "Deepanjal" run into the preceding run (leave "S." and rest of paragraph untouched).
ReplaceExact "According to a study in Korea entitled, ""Study and Evaluation of Tourism Websites based on User Perspective"" by Deepanjal " "According to a study in Korea entitled, ""Study and Evaluation of Tourism Websites based on User Perspective"" by Deepanjal "

# 3. Add a comma after "analysis" (splits into new runs naturally via COM edit).
ReplaceExact "content analysis and find it difficult" "content analysis, and find it difficult"

# 4. Merge "Palkoska" run into surrounding runs (leave "electronic" etc. untouched).
ReplaceExact "A study conducted by J. Palkoska, et al. in 2000 regarding " "A study conducted by J. Palkoska, et al. in 2000 regarding "

# 5. Rework the Overtourism sentence: move "in 2018" to the front and split "This phenomenon" -> "It".
ReplaceExact "Overtourism began to emerge as a serious and dangerous phenomenon in 2018. This phenomenon has impacted several cities, cultural heritage sites, recreational areas, and islands." "In 2018, Overtourism began to emerge as a serious and dangerous phenomenon. It has impacted several cities, cultural heritage sites, recreational areas, and islands."

# 6. Merge "Novabos" run into the following run.
ReplaceExact "Novabos et al., 2015, researchers at the University of the Philippines Diliman in Quezon City proposed a comprehensive and reliable instrument for measuring the perceived quality of destination websites." "Novabos et al., 2015, researchers at the University of the Philippines Diliman in Quezon City proposed a comprehensive and reliable instrument for measuring the perceived quality of destination websites."

# 7. Merge "Lehto" run into the following run.
ReplaceExact "Lehto, X.Y., Kim, D.Y. & Morrison, A.M. (2006), found that 93% of Internet users who seek travel information online visit official tourism websites." "Lehto, X.Y., Kim, D.Y. & Morrison, A.M. (2006), found that 93% of Internet users who seek travel information online visit official tourism websites."
